$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "2. Data reporter" section (B6:B10) with the new organization's
# contact information.
$ws.Range("B6").Value  = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value  = "Kalymbetova Yryskan"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Reflect the updated selection left by the author (B6:B10, active cell B6).
$ws.Range("B6:B10").Select()
